# Auto-generated edit script: applies market-data value updates to the
# "Chocobo_Profits" workbook (8 job sheets: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# For each changed row, columns H-N hold market/profit figures that the
# scheduled runner refreshed. A couple of rows also gained or lost a trailing
# cell (LeveProfitNQ/HQ only exists where it is meaningful).

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(55, 8).Value = 384.54544
$ws.Cells.Item(55, 9).Value = 257.18182
$ws.Cells.Item(55, 10).Value = 511.9091
$ws.Cells.Item(55, 11).Value = 257.18182
$ws.Cells.Item(55, 12).Value = 511.9091
$ws.Cells.Item(55, 13).Value = -43.18182000000002
$ws.Cells.Item(55, 14).Value = -939.9091000000001
$ws.Cells.Item(94, 8).Value = 2895
$ws.Cells.Item(94, 9).Value = 2895
$ws.Cells.Item(94, 11).Value = 2895
$ws.Cells.Item(94, 13).Value = -2444
$ws.Cells.Item(106, 8).Value = 1790.25
$ws.Cells.Item(106, 9).Value = 1416.4783
$ws.Cells.Item(106, 11).Value = 1416.4783
$ws.Cells.Item(106, 13).Value = -785.4783
$ws.Cells.Item(107, 8).Value = 1049.08
$ws.Cells.Item(107, 9).Value = 1049.08
$ws.Cells.Item(107, 10).Value = 0
$ws.Cells.Item(107, 11).Value = 1049.08
$ws.Cells.Item(107, 12).Value = 0
$ws.Cells.Item(107, 13).Value = 870.9200000000001
$ws.Cells.Item(107, 14).ClearContents()
$ws.Cells.Item(112, 8).Value = 1384.9155
$ws.Cells.Item(112, 10).Value = 1384.9155
$ws.Cells.Item(112, 12).Value = 4154.7465
$ws.Cells.Item(112, 14).Value = -6370.7465
$ws.Cells.Item(123, 8).Value = 41890
$ws.Cells.Item(123, 10).Value = 41890
$ws.Cells.Item(123, 12).Value = 41890
$ws.Cells.Item(123, 14).Value = -51690
$ws.Cells.Item(125, 8).Value = 1564
$ws.Cells.Item(125, 9).Value = 1016
$ws.Cells.Item(125, 10).Value = 1746.6666
$ws.Cells.Item(125, 11).Value = 9144
$ws.Cells.Item(125, 12).Value = 15719.9994
$ws.Cells.Item(125, 13).Value = -6684
$ws.Cells.Item(125, 14).Value = -20639.9994
$ws.Cells.Item(132, 8).Value = 26421272
$ws.Cells.Item(132, 9).Value = 30422918
$ws.Cells.Item(132, 11).Value = 91268754
$ws.Cells.Item(132, 13).Value = -91266224
$ws.Cells.Item(138, 8).Value = 1520.0834
$ws.Cells.Item(138, 9).Value = 1124.1
$ws.Cells.Item(138, 10).Value = 3500
$ws.Cells.Item(138, 11).Value = 3372.3
$ws.Cells.Item(138, 12).Value = 10500
$ws.Cells.Item(138, 13).Value = 1767.7
$ws.Cells.Item(138, 14).Value = -20780

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 3085.1
$ws.Cells.Item(132, 9).Value = 1592
$ws.Cells.Item(132, 10).Value = 5324.75
$ws.Cells.Item(132, 11).Value = 4776
$ws.Cells.Item(132, 12).Value = 15974.25
$ws.Cells.Item(132, 13).Value = -2246
$ws.Cells.Item(132, 14).Value = -21034.25
$ws.Cells.Item(137, 8).Value = 43262.5
$ws.Cells.Item(137, 10).Value = 43262.5
$ws.Cells.Item(137, 12).Value = 43262.5
$ws.Cells.Item(137, 14).Value = -53462.5

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(59, 8).Value = 59000
$ws.Cells.Item(59, 10).Value = 59000
$ws.Cells.Item(59, 12).Value = 59000
$ws.Cells.Item(59, 14).Value = -60694
$ws.Cells.Item(94, 8).Value = 436.4
$ws.Cells.Item(94, 9).Value = 436.4
$ws.Cells.Item(94, 10).Value = 0
$ws.Cells.Item(94, 11).Value = 436.4
$ws.Cells.Item(94, 12).Value = 0
$ws.Cells.Item(94, 13).Value = 14.60000000000002
$ws.Cells.Item(94, 14).ClearContents()
$ws.Cells.Item(99, 8).Value = 1000.3889
$ws.Cells.Item(99, 9).Value = 894.5333000000001
$ws.Cells.Item(99, 10).Value = 1529.6666
$ws.Cells.Item(99, 11).Value = 894.5333000000001
$ws.Cells.Item(99, 12).Value = 1529.6666
$ws.Cells.Item(99, 13).Value = 603.4666999999999
$ws.Cells.Item(99, 14).Value = -4525.6666
$ws.Cells.Item(105, 8).Value = 1795.2028
$ws.Cells.Item(105, 9).Value = 1671.3539
$ws.Cells.Item(105, 10).Value = 2689.6667
$ws.Cells.Item(105, 11).Value = 1671.3539
$ws.Cells.Item(105, 12).Value = 2689.6667
$ws.Cells.Item(105, 13).Value = 75.64609999999993
$ws.Cells.Item(105, 14).Value = -6183.6667
$ws.Cells.Item(107, 8).Value = 1299.75
$ws.Cells.Item(107, 9).Value = 1240.1177
$ws.Cells.Item(107, 11).Value = 1240.1177
$ws.Cells.Item(107, 13).Value = 679.8823
$ws.Cells.Item(134, 8).Value = 4549.294
$ws.Cells.Item(134, 9).Value = 1945
$ws.Cells.Item(134, 10).Value = 5350.615
$ws.Cells.Item(134, 11).Value = 5835
$ws.Cells.Item(134, 12).Value = 16051.845
$ws.Cells.Item(134, 13).Value = -3300
$ws.Cells.Item(134, 14).Value = -21121.845
$ws.Cells.Item(137, 8).Value = 50773.6
$ws.Cells.Item(137, 10).Value = 50773.6
$ws.Cells.Item(137, 12).Value = 50773.6
$ws.Cells.Item(137, 14).Value = -60973.6

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 170828.53
$ws.Cells.Item(31, 9).Value = 347364.3
$ws.Cells.Item(31, 10).Value = 2904.244
$ws.Cells.Item(31, 11).Value = 347364.3
$ws.Cells.Item(31, 12).Value = 2904.244
$ws.Cells.Item(31, 13).Value = -347069.3
$ws.Cells.Item(31, 14).Value = -3494.244
$ws.Cells.Item(34, 8).Value = 170828.53
$ws.Cells.Item(34, 9).Value = 347364.3
$ws.Cells.Item(34, 10).Value = 2904.244
$ws.Cells.Item(34, 11).Value = 347364.3
$ws.Cells.Item(34, 12).Value = 2904.244
$ws.Cells.Item(34, 13).Value = -347162.3
$ws.Cells.Item(34, 14).Value = -3308.244
$ws.Cells.Item(105, 8).Value = 2357
$ws.Cells.Item(105, 9).Value = 2627.25
$ws.Cells.Item(105, 10).Value = 1996.6666
$ws.Cells.Item(105, 11).Value = 2627.25
$ws.Cells.Item(105, 12).Value = 1996.6666
$ws.Cells.Item(105, 13).Value = -880.25
$ws.Cells.Item(105, 14).Value = -5490.6666
$ws.Cells.Item(107, 8).Value = 902.1177
$ws.Cells.Item(107, 9).Value = 472.35294
$ws.Cells.Item(107, 10).Value = 1331.8823
$ws.Cells.Item(107, 11).Value = 472.35294
$ws.Cells.Item(107, 12).Value = 1331.8823
$ws.Cells.Item(107, 13).Value = 1447.64706
$ws.Cells.Item(107, 14).Value = -5171.8823

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value = 3449.4666
$ws.Cells.Item(68, 9).Value = 1078.25
$ws.Cells.Item(68, 10).Value = 5346.44
$ws.Cells.Item(68, 11).Value = 3234.75
$ws.Cells.Item(68, 12).Value = 16039.32
$ws.Cells.Item(68, 13).Value = -2423.75
$ws.Cells.Item(68, 14).Value = -17661.32
$ws.Cells.Item(71, 8).Value = 3449.4666
$ws.Cells.Item(71, 9).Value = 1078.25
$ws.Cells.Item(71, 10).Value = 5346.44
$ws.Cells.Item(71, 11).Value = 9704.25
$ws.Cells.Item(71, 12).Value = 48117.96
$ws.Cells.Item(71, 13).Value = -5648.25
$ws.Cells.Item(71, 14).Value = -56229.96
$ws.Cells.Item(107, 8).Value = 9824796
$ws.Cells.Item(107, 9).Value = 349.79166
$ws.Cells.Item(107, 11).Value = 1049.37498
$ws.Cells.Item(107, 13).Value = 870.6250199999999
$ws.Cells.Item(113, 8).Value = 465.06818
$ws.Cells.Item(113, 9).Value = 471.18518
$ws.Cells.Item(113, 10).Value = 455.35294
$ws.Cells.Item(113, 11).Value = 1413.55554
$ws.Cells.Item(113, 12).Value = 1366.05882
$ws.Cells.Item(113, 13).Value = 756.4444599999999
$ws.Cells.Item(113, 14).Value = -5706.05882
$ws.Cells.Item(131, 8).Value = 5000851.5
$ws.Cells.Item(131, 10).Value = 920.41174
$ws.Cells.Item(131, 12).Value = 2761.23522
$ws.Cells.Item(131, 14).Value = -12841.23522

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(46, 8).Value = 23997
$ws.Cells.Item(46, 10).Value = 24747.166
$ws.Cells.Item(46, 12).Value = 24747.166
$ws.Cells.Item(46, 14).Value = -25059.166
$ws.Cells.Item(113, 8).Value = 1329.3334
$ws.Cells.Item(113, 9).Value = 1300
$ws.Cells.Item(113, 10).Value = 1344
$ws.Cells.Item(113, 11).Value = 1300
$ws.Cells.Item(113, 12).Value = 1344
$ws.Cells.Item(113, 13).Value = 870
$ws.Cells.Item(113, 14).Value = -5684
$ws.Cells.Item(126, 8).Value = 3468.9368
$ws.Cells.Item(126, 9).Value = 2808.423
$ws.Cells.Item(126, 10).Value = 4741.037
$ws.Cells.Item(126, 11).Value = 8425.269
$ws.Cells.Item(126, 12).Value = 14223.111
$ws.Cells.Item(126, 13).Value = -5955.269
$ws.Cells.Item(126, 14).Value = -19163.111
$ws.Cells.Item(137, 8).Value = 39073.332
$ws.Cells.Item(137, 10).Value = 48610
$ws.Cells.Item(137, 12).Value = 48610
$ws.Cells.Item(137, 14).Value = -58810

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(118, 8).Value = 24769.166
$ws.Cells.Item(118, 10).Value = 24769.166
$ws.Cells.Item(118, 12).Value = 24769.166
$ws.Cells.Item(118, 14).Value = -28083.166

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 545.85
$ws.Cells.Item(107, 9).Value = 378.66666
$ws.Cells.Item(107, 10).Value = 893.0769
$ws.Cells.Item(107, 11).Value = 1135.99998
$ws.Cells.Item(107, 12).Value = 2679.2307
$ws.Cells.Item(107, 13).Value = 784.0000199999999
$ws.Cells.Item(107, 14).Value = -6519.2307
$ws.Cells.Item(116, 8).Value = 35000
$ws.Cells.Item(116, 10).Value = 35000
$ws.Cells.Item(116, 12).Value = 35000
$ws.Cells.Item(116, 14).Value = -44178
$ws.Cells.Item(122, 8).Value = 7388.8887
$ws.Cells.Item(122, 9).Value = 2900
$ws.Cells.Item(122, 10).Value = 7950
$ws.Cells.Item(122, 11).Value = 8700
$ws.Cells.Item(122, 12).Value = 23850
$ws.Cells.Item(122, 14).Value = -28750
$ws.Cells.Item(122, 13).Value = -6250
$ws.Cells.Item(136, 8).Value = 3271.44
$ws.Cells.Item(136, 9).Value = 643.3889
$ws.Cells.Item(136, 11).Value = 1930.1667
$ws.Cells.Item(136, 13).Value = 619.8332999999998

Write-Output "Applied 201 value updates, 1 additions, 2 removals."
